$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Value2 -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
}
